$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 400 (shifts existing rows 400..495 down to 401..496)
$ws.Rows.Item(400).Insert()

# Populate the newly inserted row 400 with the new data record
$ws.Range("A400").Value = 10
$ws.Range("B400").Value = "Vega Modelo de Temuco"
$ws.Range("C400").Value = "La Araucanía"
$ws.Range("D400").Value = 44932
$ws.Range("E400").Value = 9
$ws.Range("F400").Value = 100112040
$ws.Range("G400").Value = "Cilantro"
$ws.Range("H400").Value = "Sin especificar"
$ws.Range("I400").Value = "Primera"
$ws.Range("J400").Value = 20
$ws.Range("K400").Value = 10000
$ws.Range("L400").Value = 10000
$ws.Range("M400").Value = 10000
$ws.Range("N400").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O400").Value = "Provincia de Cautín"
$ws.Range("P400").Value = 5000
$ws.Range("Q400").Value = 2
$ws.Range("R400").Value = "Hortaliza"

# Ensure the date cell keeps the date/time number format used by the rest of column D
$ws.Range("D400").NumberFormat = "YYYY-MM-DD HH:MM:SS"
